# Update profit.py after running on 2025-09-02
#
# 1. Drop the (now unused) "Sheet2" scratch sheet.
# 2. Append the new day's row to Sheet1: 09/02/2025 -> 13336.21
$wb = $excel.ActiveWorkbook

[void]$wb.Worksheets("Sheet2").Delete()

$ws = $wb.Worksheets("Sheet1")

# Re-apply A5's existing date/time format so the style table collapses its
# duplicate entry (the style index this cell points at shifts down by one).
$ws.Range("A5").NumberFormat = $ws.Range("A5").NumberFormat

$nextRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Keep the date as plain text (matching the other rows in column A)
# instead of letting Excel auto-convert "09/02/2025" into a date serial.
$dateCell = $ws.Cells.Item($nextRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/02/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($nextRow, 2).Value = 13336.21
